$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.893.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.629.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E8').Value = '  -0.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.258'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0608'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.860.56'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.612.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.904.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.97'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.07'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.42'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.413.22'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('E35').Value = '  +2.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.10%  '
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.769.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('E46').Value = '  -3.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.101'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0503'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.62'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.26%  '
